# epexspot_prices.xlsx — daily auto-update
# 1) "Prix Spot": a new day ("08-nov") of spot prices arrived. Its column is
#    inserted right before the "01-oct." column (DK), shifting every later
#    date column one slot to the right. The new column has no price data yet
#    ("-") for every hour row.
# 2) "Gaz" / "CO2": append the next day's closing price as a new last row.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Prix Spot" ------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before DK (shifts DK:EO -> DL:EP).
$ws.Range("DK1").EntireColumn.Insert()

# New header cell picks up the header style automatically; give it the date.
$ws.Range("DK1").Value = "08-nov"

# No data yet for this new day on any of the 24 hour rows.
$ws.Range("DK2:DK25").Value = "-"

# ---- Sheet "Gaz" -------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Keep the date as literal text (matches every other row in the column)
# instead of letting it be auto-recognised as a date value.
$wsGaz.Range("A144").NumberFormat = "@"
$wsGaz.Range("A144").Value = "2025-11-06"
$wsGaz.Range("B144").Value = 30.35

# ---- Sheet "CO2" --------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A144").NumberFormat = "@"
$wsCO2.Range("A144").Value = "2025-11-06"
$wsCO2.Range("B144").Value = 79.94
